# Updates cryptos list values (price/volume columns) to match the latest
# scraped data, as produced by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.929.03'
$ws.Range('E2').Value = '  -1.82%  '
$ws.Range('D3').Value = '1.867.32'
$ws.Range('E3').Value = '  -2.58%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.90'
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4994'
$ws.Range('E7').Value = '  -2.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3827'
$ws.Range('E8').Value = '  -3.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08875'
$ws.Range('E9').Value = '  -8.76%  '
$ws.Range('E10').Value = '  -2.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.54'
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.370'
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.63'
$ws.Range('E13').Value = '  -1.93%  '
$ws.Range('D14').Value = '1.865.42'
$ws.Range('E14').Value = '  -3.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.229'
$ws.Range('E15').Value = '  -2.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001096'
$ws.Range('E17').Value = '  -3.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '91.13'
$ws.Range('E18').Value = '  -2.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06665'
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.97'
$ws.Range('E20').Value = '  -0.87%  '
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.109'
$ws.Range('E22').Value = '  -2.33%  '
$ws.Range('D23').Value = '27.917.64'
$ws.Range('E23').Value = '  -2.02%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.285'
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('D26').Value = '2.078.33'
$ws.Range('E26').Value = '  -3.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.501'
$ws.Range('E27').Value = '  -6.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.82'
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.66'
$ws.Range('E29').Value = '  -2.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.42'
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1058'
$ws.Range('E31').Value = '  -1.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.054'
$ws.Range('E32').Value = '  -4.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.591'
$ws.Range('E33').Value = '  -1.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.606'
$ws.Range('E34').Value = '  -0.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.460'
$ws.Range('E35').Value = '  -3.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06526'
$ws.Range('E36').Value = '  -2.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02387'
$ws.Range('E37').Value = '  -2.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2180'
$ws.Range('E38').Value = '  -1.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.278'
$ws.Range('E39').Value = '  +4.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.200'
$ws.Range('E40').Value = '  -4.51%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.45'
$ws.Range('E41').Value = '  -1.52%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6341'
$ws.Range('E42').Value = '  -1.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.895'
$ws.Range('E43').Value = '  -3.18%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.17'
$ws.Range('E45').Value = '  -3.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5990'
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.280'
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.670'
$ws.Range('E48').Value = '  -2.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.227'
$ws.Range('E49').Value = '  +2.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.979'
$ws.Range('E50').Value = '  -3.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.93'
$ws.Range('E51').Value = '  -3.03%  '
